$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" (strikeouts -> K) values for column G, rows 2-16
$values = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 3
    15 = 1
    16 = 0
    17 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
